$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C4").Value = -0.02021663149666526
$ws.Range("E4").Value = 0.00944277877771115
$ws.Range("F4").Value = 0.01953305233332209
$ws.Range("H4").Value = -0.01932764957310598
$ws.Range("J4").Value = -0.005341641227229319

$ws.Range("C5").Value = -0.009474706458988257
$ws.Range("E5").Value = -0.006712324012492959
$ws.Range("F5").Value = 0.002650902442036097
$ws.Range("H5").Value = 0.001057084842283393
$ws.Range("J5").Value = 0.01273700181140178

$ws.Range("C6").Value = 0.005435903641436144
$ws.Range("E6").Value = -0.0004468193458727738
$ws.Range("F6").Value = 0.003988895679555828
$ws.Range("H6").Value = -0.004334956877398274
$ws.Range("J6").Value = 0.006552475632849239

$ws.Range("C7").Value = 0.003723509524940381
$ws.Range("E7").Value = -0.02004933142597325
$ws.Range("F7").Value = 0.0005657480866299233
$ws.Range("H7").Value = 0.005263905714556228
$ws.Range("J7").Value = 0.01445295345237692

$ws.Range("C8").Value = 0.07802736340909452
$ws.Range("E8").Value = -0.2246942173077686
$ws.Range("F8").Value = -0.9629001014600039
$ws.Range("H8").Value = 0.9999999080319962
$ws.Range("J8").Value = -0.04040152520516974

$ws.Range("C9").Value = 0.9709072535562899
$ws.Range("E9").Value = 0.01628305361132214
$ws.Range("F9").Value = 0.02394338918173557
$ws.Range("H9").Value = -0.02724055961762238
$ws.Range("J9").Value = 0.01638948656717849

$ws.Range("C10").Value = 0.005281818835272753
$ws.Range("E10").Value = -0.01527089561883582
$ws.Range("F10").Value = -0.006660075914403036
$ws.Range("H10").Value = 0.01121677052867082
$ws.Range("J10").Value = -0.01770858189635154

$ws.Range("C11").Value = -0.002057653906306156
$ws.Range("E11").Value = 0.008232181193287245
$ws.Range("F11").Value = -0.007694576659783065
$ws.Range("H11").Value = 0.005148021805920871
$ws.Range("J11").Value = -0.004990855374187153

$ws.Range("C12").Value = 0.06920738052829521
$ws.Range("E12").Value = 0.0001128647085145883
$ws.Range("F12").Value = -0.06565128137805123
$ws.Range("H12").Value = 0.02383454898538195
$ws.Range("J12").Value = 0.01796177608109947

$ws.Range("C13").Value = 0.06714553881382154
$ws.Range("E13").Value = 0.8516846920193875
$ws.Range("F13").Value = -0.245630975553239
$ws.Range("H13").Value = -0.01319457604778304
$ws.Range("J13").Value = -0.00169264213521143

$ws.Range("C14").Value = -0.163568025614721
$ws.Range("E14").Value = 0.02252665213306608
$ws.Range("F14").Value = -0.02517822091112883
$ws.Range("H14").Value = 0.0207201375968055
$ws.Range("J14").Value = -0.0005100817339697825

$ws.Range("C15").Value = -0.01453519670940787
$ws.Range("E15").Value = 0.008239000457560017
$ws.Range("F15").Value = -0.003166768542670742
$ws.Range("H15").Value = 0.002036638353465534
$ws.Range("J15").Value = 0.00255185825175936

$ws.Range("C16").Value = -0.005309955860398234
$ws.Range("E16").Value = -0.02103744487349779
$ws.Range("F16").Value = -0.02111889665275586
$ws.Range("H16").Value = 0.02683284952131398
$ws.Range("J16").Value = -0.01377416873392978

$ws.Range("C17").Value = 0.006002090448083616
$ws.Range("E17").Value = 0.01878996132759845
$ws.Range("F17").Value = -0.04778840197553607
$ws.Range("H17").Value = 0.0411947475997899
$ws.Range("J17").Value = 0.002587143452201412

$ws.Range("C18").Value = 0.02663246429729857
$ws.Range("E18").Value = 0.0156903449796138
$ws.Range("F18").Value = -0.01690512134820485
$ws.Range("H18").Value = 0.01337277960691118
$ws.Range("J18").Value = 0.02388859927857544

$ws.Range("C19").Value = 0.01085084203403368
$ws.Range("E19").Value = 0.0005845856873834275
$ws.Range("F19").Value = -0.003659187218367488
$ws.Range("H19").Value = 0.002063762866550514
$ws.Range("J19").Value = 0.004227246625160694

$ws.Range("C20").Value = 0.009232484049299361
$ws.Range("E20").Value = -0.005458068314322732
$ws.Range("F20").Value = 0.001656151170246047
$ws.Range("H20").Value = -0.0005838196073527842
$ws.Range("J20").Value = -0.007931675869542396

$ws.Range("C21").Value = 0.02367959374718374
$ws.Range("E21").Value = 0.006257330458293217
$ws.Range("F21").Value = -0.02467313042692521
$ws.Range("H21").Value = 0.02410244486809779
$ws.Range("J21").Value = -0.01621114380055195

$ws.Range("C22").Value = 0.01217714621508585
$ws.Range("E22").Value = -0.001810838952433558
$ws.Range("F22").Value = 0.001723502660940106
$ws.Range("H22").Value = -0.0007684522867380914
$ws.Range("J22").Value = -0.005064764741794004

$ws.Range("C23").Value = -0.008212455880498234
$ws.Range("E23").Value = 0.00278359076734363
$ws.Range("F23").Value = 0.009292844627713783
$ws.Range("H23").Value = -0.008494314963772598
$ws.Range("J23").Value = 0.01352005997434009
